$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Occurrence 1: Professional summary paragraph.
# Simple in-run text substitution - no formatting change needed.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters,",
    2
)

# ---------------------------------------------------------------------
# Occurrence 2: "Partner - Siege Analytics" bullet point.
# "50M" needs to be its own bold, colored run (matching the styling
# already used for the "23%"/"64%" figures in the same sentence), so we
# locate the phrase via paragraph text offsets and use a Range so Word
# naturally splits the run where the formatting changes.
# ---------------------------------------------------------------------
$targetParaText = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Discovered systematic race coding errors*") {
        $targetParaText = $p
        break
    }
}

$pStart = $targetParaText.Range.Start
$fullText = $targetParaText.Range.Text
$oldPhrase = "all Black and Asian-American voters"
$newPhrase = "50M voters"
$boldWord = "50M"

$idx = $fullText.IndexOf($oldPhrase)
$rangeStart = $pStart + $idx
$rangeEnd = $rangeStart + $oldPhrase.Length

# Replace the phrase text first (keeps a single run for now).
$sub = $d.Range($rangeStart, $rangeEnd)
$sub.Text = $newPhrase

# Now re-scope just the "50M" portion and give it bold + the accent color
# used elsewhere in this document for highlighted statistics (2C3E50).
$boldStart = $rangeStart
$boldEnd = $rangeStart + $boldWord.Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = 1
$boldRange.Font.Color = 5258796

# ---------------------------------------------------------------------
# Occurrence 3: Project impact statement.
# Simple in-run text substitution - no formatting change needed.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2
)

Write-Output "edits applied"
